# Applies:
#  1. "Figure out .gitignore for project" -> split into runs with
#     spell-check proofErr markers bracketing "gitignore".
#  2. "Add empty space lines to the output (check cmd prompt output, ...)"
#     -> paragraph (and every run) struck through, plus proofErr markers
#     bracketing "cmd".
#  3. "Display team when done" -> paragraph (and its run) struck through.

$d = $word.ActiveDocument

function Find-ParagraphByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

function Get-ParaOpenTag($paragraph) {
    # Pull the paragraph's own <w:p ...> opening tag (with its original
    # w14:paraId/w:rsidR/etc. attributes intact) out of its WordOpenXML,
    # so the replacement we splice back in keeps the same identity.
    $oxml = $paragraph.Range.WordOpenXML
    if ($oxml -match '<w:body>(?<tag><w:p\b[^>]*>)') {
        return $matches['tag']
    }
    return '<w:p>'
}

function Set-ParagraphXml($paragraph, $openTag, $innerXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
          '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:xml="http://www.w3.org/XML/1998/namespace">' +
              '<w:body>' + $openTag + $innerXml + '</w:p></w:body>' +
            '</w:document>' +
          '</pkg:xmlData>' +
        '</pkg:part>' +
      '</pkg:package>'
    $paragraph.Range.InsertXML($pkg)
}

# --- 1. "Figure out .gitignore for project" ---------------------------
$p1 = Find-ParagraphByText("Figure out .gitignore for project")
$tag1 = Get-ParaOpenTag($p1)
$inner1 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
  '<w:r><w:t>Figure out .</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>gitignore</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> for project</w:t></w:r>'
Set-ParagraphXml $p1 $tag1 $inner1

# --- 2. "Add empty space lines to the output ..." ----------------------
$p2 = Find-ParagraphByText("Add empty space lines to the output")
$tag2 = Get-ParaOpenTag($p2)
$inner2 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Add empty space lines to the output (check </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:strike/></w:rPr><w:t>cmd</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> prompt output, right now it’s jumbled and a little hard to read, just add some empty lines to clean up the output)</w:t></w:r>'
Set-ParagraphXml $p2 $tag2 $inner2

# --- 3. "Display team when done" ---------------------------------------
$p3 = Find-ParagraphByText("Display team when done")
$tag3 = Get-ParaOpenTag($p3)
$inner3 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:strike/></w:rPr><w:t>Display team when done</w:t></w:r>'
Set-ParagraphXml $p3 $tag3 $inner3

Write-Output "done"
